# Add a new "2022-Q3" sheet (cloned from "2021-Q4" so it inherits the same
# header/style layout), positioned right after "总计" and before "2021-Q4".
# Existing "2021-Q4", "2021-Q2", "2021-Q1" sheets are left untouched and
# simply shift one position to the right.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($q4)

# The duplicate is inserted immediately before the original (which keeps its
# own name/position) and becomes the active sheet; rename it to the new
# quarter. It sits right after "总计" ⇒ worksheet index 2.
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# --- Populate the new "2022-Q3" sheet with the fund data for the quarter ---
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).NumberFormat = "@"
$q3.Cells.Item(2,2).Value = "014062"
$q3.Cells.Item(2,3).NumberFormat = "@"
$q3.Cells.Item(2,3).Value = "景顺长城专精特新量化优选股票A"
$q3.Cells.Item(2,4).NumberFormat = "@"
$q3.Cells.Item(2,4).Value = "8.02"
$q3.Cells.Item(2,5).NumberFormat = "@"
$q3.Cells.Item(2,5).Value = "91.10"
$q3.Cells.Item(2,6).NumberFormat = "@"
$q3.Cells.Item(2,6).Value = "1.66"
$q3.Cells.Item(2,7).NumberFormat = "@"
$q3.Cells.Item(2,7).Value = "0.1331"
$q3.Cells.Item(2,8).Value = 10

$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).NumberFormat = "@"
$q3.Cells.Item(3,2).Value = "014063"
$q3.Cells.Item(3,3).NumberFormat = "@"
$q3.Cells.Item(3,3).Value = "景顺长城专精特新量化优选股票C"
$q3.Cells.Item(3,4).NumberFormat = "@"
$q3.Cells.Item(3,4).Value = "5.41"
$q3.Cells.Item(3,5).NumberFormat = "@"
$q3.Cells.Item(3,5).Value = "91.10"
$q3.Cells.Item(3,6).NumberFormat = "@"
$q3.Cells.Item(3,6).Value = "1.66"
$q3.Cells.Item(3,7).NumberFormat = "@"
$q3.Cells.Item(3,7).Value = "0.0898"
$q3.Cells.Item(3,8).Value = 10

# --- Update the "总计" (summary) sheet: a new row for 2022-Q3 is inserted
# at the top of the breakdown, pushing the older quarters down by one row.
$total = $wb.Worksheets.Item("总计")

# Copy the formatting of the last existing data row down onto the new row 5
# before touching any values, so row 5 picks up the same bold/border style
# used by rows 2-4.
$total.Cells.Item(4,1).Copy()
$total.Cells.Item(5,1).PasteSpecial(-4122)

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q1"
$total.Cells.Item(5,3).Value = 2
$total.Cells.Item(5,4).Value = 0.04

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q2"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.04

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 2
$total.Cells.Item(3,4).Value = 0.09

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.22

$total.Activate()
